# VerveStacks ITA model update - 2025-09-01 16:30
#
# On the "Misc" sheet, the generic wind entries in the ~TFM_AVA and
# ~TFM_TOPINS blocks are split into separate onshore/offshore entries:
#   - ~TFM_AVA:    ELC_won* (wind-onshore availability mask) -> ELC_wo* (covers both)
#   - ~TFM_TOPINS: E[_]W*/wind (generic wind) -> E[_]WOF*/windoff (offshore),
#                  with a new E[_]WON*/windon (onshore) row added right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new onshore-wind row: inserting at row 44 pushes the
# existing ~TFM_INS block (old rows 46-49) down to rows 47-50, preserving
# their formatting (e.g. the quote-prefixed "-pasti"/"-life" cells).
$ws.Rows(44).Insert()

# ~TFM_TOPINS block: row 43 now describes offshore wind specifically.
$ws.Range("C43").Value2 = "E[_]WOF*"
$ws.Range("D43").Value2 = "windoff"

# New row 44: onshore wind counterpart.
$ws.Range("C44").Value2 = "E[_]WON*"
$ws.Range("D44").Value2 = "windon"
$ws.Range("E44").Value2 = "IN"

# ~TFM_AVA block: widen the onshore-only wildcard to match both wind units.
$ws.Range("D41").Value2 = "ELC_wo*"

# Leave the selection where the author left it when saving.
$ws.Range("D42").Select()
